$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unprotect so we can insert rows / edit data -------------------------
$ws.Unprotect("EF56")

# --- Insert 2 new rows above row 6 (pushes the old "Story 1" task rows
#     and the Ideal/Actual total rows down by two) ------------------------
$ws.Rows("6:7").Insert()

# Copy formatting from row 8 (a still-intact "task" row) onto the two new
# blank rows so they pick up the same styles (s="3" labels, s="4" numbers)
$ws.Range("A8:S8").Copy()
$ws.Range("A6:S7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Replace the logbook entries for rows 6-10 with the new "Android"
#     tasks (this also makes the old "Story 1"/"Task 1.1-1.3" shared
#     strings unreferenced so they drop out of sharedStrings.xml on save)
$ws.Range("B6").Value2 = "Android"
$ws.Range("C6").Value2 = "Pengenalan Android"
$ws.Range("D6").Value2 = "Apa itut Android ?"
$ws.Range("E6").Value2 = 2
$ws.Range("F6").Value2 = 4
$ws.Range("G6").Value2 = 0
$ws.Range("H6").Value2 = 0
$ws.Range("I6").Value2 = 0
$ws.Range("J6").Value2 = 0

$ws.Range("B7").Value2 = "Android"
$ws.Range("C7").Value2 = "Pengenalan Android"
$ws.Range("D7").Value2 = "Sedikit Sejarah Android"
$ws.Range("E7").Value2 = 2
$ws.Range("F7").Value2 = 2
$ws.Range("G7").Value2 = 3
$ws.Range("H7").Value2 = 5
$ws.Range("I7").Value2 = 2
$ws.Range("J7").Value2 = 0

$ws.Range("B8").Value2 = "Android"
$ws.Range("C8").Value2 = "Pengenalan Android"
$ws.Range("D8").Value2 = "Jelaskan ICS"
$ws.Range("E8").Value2 = 2
$ws.Range("F8:J8").ClearContents()

$ws.Range("B9").Value2 = "Android"
$ws.Range("C9").Value2 = "Pengenalan Android"
$ws.Range("D9").Value2 = "Jelaskan JB"
$ws.Range("E9").Value2 = 2
$ws.Range("F9:J9").ClearContents()

$ws.Range("B10").Value2 = "Android"
$ws.Range("C10").Value2 = "Pengenalan Android"
$ws.Range("D10").Value2 = "Jelaskan Kitkat"
$ws.Range("E10").Value2 = 2
$ws.Range("F10").Value2 = 5
$ws.Range("G10").Value2 = 3
$ws.Range("H10").Value2 = 2
$ws.Range("I10").Value2 = 2
$ws.Range("J10").Value2 = 0

# Rows 11-20 already hold the right Story2/3/4 + Task2.x-4.x data because
# they were simply shifted down two rows by the Insert() above - nothing
# further to change there.

# --- Fix up the two summary rows (now at 21 "Ideal" / 22 "Actual") so the
#     SUM()/ideal-burndown ranges cover the full new logbook (rows 6-20)
$ws.Range("E21").Formula = "=SUM(E6:E20)"
$ws.Range("F21").Formula = "=E21-`$E`$21/10"
$ws.Range("G21").Formula = "=F21-`$E`$21/10"
$ws.Range("H21").Formula = "=G21-`$E`$21/10"
$ws.Range("I21").Formula = "=H21-`$E`$21/10"
$ws.Range("J21").Formula = "=I21-`$E`$21/10"
$ws.Range("K21").Formula = "=J21-`$E`$21/10"
$ws.Range("L21").Formula = "=K21-`$E`$21/10"
$ws.Range("M21").Formula = "=L21-`$E`$21/10"
$ws.Range("N21").Formula = "=M21-`$E`$21/10"
$ws.Range("O21").Formula = "=N21-`$E`$21/10"

$ws.Range("E22").Formula = "=SUM(E6:E20)"
$ws.Range("F22").Formula = "=SUM(F6:F20)"
$ws.Range("G22").Formula = "=SUM(G6:G20)"
$ws.Range("H22").Formula = "=SUM(H6:H20)"
$ws.Range("I22").Formula = "=SUM(I6:I20)"
$ws.Range("J22").Formula = "=SUM(J6:J20)"

# --- Column width tweaks (C/D got wider) ---------------------------------
$ws.Columns("C").ColumnWidth = 20.140625
$ws.Columns("D").ColumnWidth = 23.42578125

# --- Update the chart's series formulas so they point at the relocated
#     Ideal/Actual rows (21/22 instead of 19/20) --------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$s1 = $chart.SeriesCollection().Item(1)
$s2 = $chart.SeriesCollection().Item(2)
$s1.Formula = '=SERIES("Ideal burndown",Sheet1!$F$5:$O$5,Sheet1!$F$21:$O$21,1)'
$s2.Formula = '=SERIES("Actual burndown",Sheet1!$F$5:$O$5,Sheet1!$F$22:$O$22,2)'

# --- Sheet view: scroll position / active selection moved ---------------
$ws.Range("B9").Select()
$av = $excel.ActiveWindow
$av.ScrollRow = 3

# --- Re-protect the sheet the way it was before --------------------------
$ws.Protect("EF56", $true, $true, $true, $true, $true, $true, $true, $true, $true, $false, $true, $false, $false, $false, $false)
